# Updated main GSC export data:
# - drop the oldest date row (2025-10-24) from the "Chart" sheet, shifting
#   every subsequent day's row up by one
# - the newest day (2026-01-20), which previously had no Impressions value
#   yet, now has its final value recorded

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the row for the oldest date (2025-10-24); everything below shifts up.
$ws.Rows.Item(2).Delete()

# The newest date's Impressions value is now known - fill it in on what is
# now the last row of the table.
$ws.Range("D89").Value = 79
